$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Ativação date update
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Row 10: Objetivos (PT) text update
$ws.Range("B10").Value = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."
$ws.Range("C10").Value = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."

# Row 11: Objectives (EN) new text
$ws.Range("B11").Value = "Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base"
$ws.Range("C11").Value = "Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base"

# Row 14: Programa resumido (PT) update
$ws.Range("B14").Value = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."
$ws.Range("C14").Value = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."

# Row 15: Short syllabus (EN) new text
$ws.Range("B15").Value = "Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation."
$ws.Range("C15").Value = "Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation."

# Row 16: Programa (PT) update
$ws.Range("B16").Value = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."
$ws.Range("C16").Value = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."

# Row 17: Syllabus (EN) new text
$ws.Range("B17").Value = "1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation."
$ws.Range("C17").Value = "1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation."

# Row 19: Método update
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."

# Row 20: Critério update
$ws.Range("B20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."

# Row 21: Norma de recuperação update
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."

# Row 22: Bibliografia update
$ws.Range("B22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."
$ws.Range("C22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."
